$wb = $excel.ActiveWorkbook

$managerData = @(
    ,@(1205, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://online.r-m.co.il/login.php', 0.364)
    ,@(1206, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://system.serveandcheck.com/login.php', 1.769)
    ,@(1207, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://eu.checker-soft.com/gfk-ukraine/login.php', 1.276)
    ,@(1208, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://ipsos-russia.com/login.php', 0.86)
    ,@(1209, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://tyaonline.com/login.php', 4.372)
    ,@(1210, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://marketest.checker.co.il/login.php', 1.249)
    ,@(1211, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://ms-online.co.il/login.php', 1.248)
    ,@(1212, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://www.misonline.co.il/login.php', 1.648)
    ,@(1213, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://tovanot.checker.co.il/login.php', 1.311)
    ,@(1214, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://csc.ajis-group.co.jp/login.php', 4.131)
    ,@(1215, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://ru.checker-soft.com/profpoint-ru/login.php', 0.828)
    ,@(1216, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://www.imystery.ru/login.php', 0.885)
    ,@(1217, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://eu.checker-soft.com/testing/login.php', 0.9340000000000001)
    ,@(1218, $null, '*****', $null)
    ,@(1219, '10.14.2022 12:55 (Kyiv+Israel) 09:55 (UTC) 18:55 (Japan) 15:25 (India)', 'https://online.r-m.co.il/login.php', 0.378)
    ,@(1220, '10.14.2022 12:55 (Kyiv+Israel) 09:55 (UTC) 18:55 (Japan) 15:25 (India)', 'https://system.serveandcheck.com/login.php', 1.344)
    ,@(1221, '10.14.2022 12:55 (Kyiv+Israel) 09:55 (UTC) 18:55 (Japan) 15:25 (India)', 'https://eu.checker-soft.com/gfk-ukraine/login.php', 2.565)
    ,@(1222, '10.14.2022 12:55 (Kyiv+Israel) 09:55 (UTC) 18:55 (Japan) 15:25 (India)', 'https://ipsos-russia.com/login.php', 1.253)
    ,@(1223, '10.14.2022 12:55 (Kyiv+Israel) 09:55 (UTC) 18:55 (Japan) 15:25 (India)', 'https://tyaonline.com/login.php', 4.093)
    ,@(1224, '10.14.2022 12:55 (Kyiv+Israel) 09:55 (UTC) 18:55 (Japan) 15:25 (India)', 'https://marketest.checker.co.il/login.php', 1.544)
    ,@(1225, '10.14.2022 12:56 (Kyiv+Israel) 09:56 (UTC) 18:56 (Japan) 15:26 (India)', 'https://ms-online.co.il/login.php', 1.236)
    ,@(1226, '10.14.2022 12:56 (Kyiv+Israel) 09:56 (UTC) 18:56 (Japan) 15:26 (India)', 'https://www.misonline.co.il/login.php', 1.367)
    ,@(1227, '10.14.2022 12:56 (Kyiv+Israel) 09:56 (UTC) 18:56 (Japan) 15:26 (India)', 'https://tovanot.checker.co.il/login.php', 1.205)
    ,@(1228, '10.14.2022 12:56 (Kyiv+Israel) 09:56 (UTC) 18:56 (Japan) 15:26 (India)', 'https://csc.ajis-group.co.jp/login.php', 4.673)
    ,@(1229, '10.14.2022 12:56 (Kyiv+Israel) 09:56 (UTC) 18:56 (Japan) 15:26 (India)', 'https://ru.checker-soft.com/profpoint-ru/login.php', 1.201)
    ,@(1230, '10.14.2022 12:56 (Kyiv+Israel) 09:56 (UTC) 18:56 (Japan) 15:26 (India)', 'https://www.imystery.ru/login.php', 1.212)
    ,@(1231, '10.14.2022 12:56 (Kyiv+Israel) 09:56 (UTC) 18:56 (Japan) 15:26 (India)', 'https://eu.checker-soft.com/testing/login.php', 4.37)
    ,@(1232, $null, '*****', $null)
    ,@(1233, '10.14.2022 15:55 (Kyiv+Israel) 12:55 (UTC) 21:55 (Japan) 18:25 (India)', 'https://online.r-m.co.il/login.php', 0.32)
    ,@(1234, '10.14.2022 15:55 (Kyiv+Israel) 12:55 (UTC) 21:55 (Japan) 18:25 (India)', 'https://system.serveandcheck.com/login.php', 1.239)
    ,@(1235, '10.14.2022 15:55 (Kyiv+Israel) 12:55 (UTC) 21:55 (Japan) 18:25 (India)', 'https://eu.checker-soft.com/gfk-ukraine/login.php', 1.169)
    ,@(1236, '10.14.2022 15:55 (Kyiv+Israel) 12:55 (UTC) 21:55 (Japan) 18:25 (India)', 'https://ipsos-russia.com/login.php', 0.866)
    ,@(1237, '10.14.2022 15:56 (Kyiv+Israel) 12:56 (UTC) 21:56 (Japan) 18:26 (India)', 'https://tyaonline.com/login.php', 4.07)
    ,@(1238, '10.14.2022 15:56 (Kyiv+Israel) 12:56 (UTC) 21:56 (Japan) 18:26 (India)', 'https://marketest.checker.co.il/login.php', 1.359)
    ,@(1239, '10.14.2022 15:56 (Kyiv+Israel) 12:56 (UTC) 21:56 (Japan) 18:26 (India)', 'https://ms-online.co.il/login.php', 2.771)
    ,@(1240, '10.14.2022 15:56 (Kyiv+Israel) 12:56 (UTC) 21:56 (Japan) 18:26 (India)', 'https://www.misonline.co.il/login.php', 1.484)
    ,@(1241, '10.14.2022 15:56 (Kyiv+Israel) 12:56 (UTC) 21:56 (Japan) 18:26 (India)', 'https://tovanot.checker.co.il/login.php', 1.803)
    ,@(1242, '10.14.2022 15:56 (Kyiv+Israel) 12:56 (UTC) 21:56 (Japan) 18:26 (India)', 'https://csc.ajis-group.co.jp/login.php', 4.276)
    ,@(1243, '10.14.2022 15:56 (Kyiv+Israel) 12:56 (UTC) 21:56 (Japan) 18:26 (India)', 'https://ru.checker-soft.com/profpoint-ru/login.php', 0.999)
    ,@(1244, '10.14.2022 15:57 (Kyiv+Israel) 12:57 (UTC) 21:57 (Japan) 18:27 (India)', 'https://www.imystery.ru/login.php', 1.157)
    ,@(1245, '10.14.2022 15:57 (Kyiv+Israel) 12:57 (UTC) 21:57 (Japan) 18:27 (India)', 'https://eu.checker-soft.com/testing/login.php', 0.916)
    ,@(1246, $null, '*****', $null)
    ,@(1247, '10.14.2022 19:31 (Kyiv+Israel) 16:31 (UTC) 01:31 (Japan) 22:01 (India)', 'https://online.r-m.co.il/login.php', 0.293)
    ,@(1248, '10.14.2022 19:31 (Kyiv+Israel) 16:31 (UTC) 01:31 (Japan) 22:01 (India)', 'https://system.serveandcheck.com/login.php', 1.028)
    ,@(1249, '10.14.2022 19:31 (Kyiv+Israel) 16:31 (UTC) 01:31 (Japan) 22:01 (India)', 'https://eu.checker-soft.com/gfk-ukraine/login.php', 1.066)
    ,@(1250, '10.14.2022 19:31 (Kyiv+Israel) 16:31 (UTC) 01:31 (Japan) 22:01 (India)', 'https://ipsos-russia.com/login.php', 1.033)
    ,@(1251, '10.14.2022 19:31 (Kyiv+Israel) 16:31 (UTC) 01:31 (Japan) 22:01 (India)', 'https://tyaonline.com/login.php', 3.965)
    ,@(1252, '10.14.2022 19:31 (Kyiv+Israel) 16:31 (UTC) 01:31 (Japan) 22:01 (India)', 'https://marketest.checker.co.il/login.php', 1.817)
    ,@(1253, '10.14.2022 19:32 (Kyiv+Israel) 16:32 (UTC) 01:32 (Japan) 22:02 (India)', 'https://ms-online.co.il/login.php', 2.293)
    ,@(1254, '10.14.2022 19:32 (Kyiv+Israel) 16:32 (UTC) 01:32 (Japan) 22:02 (India)', 'https://www.misonline.co.il/login.php', 1.227)
    ,@(1255, '10.14.2022 19:32 (Kyiv+Israel) 16:32 (UTC) 01:32 (Japan) 22:02 (India)', 'https://tovanot.checker.co.il/login.php', 2.051)
    ,@(1256, '10.14.2022 19:32 (Kyiv+Israel) 16:32 (UTC) 01:32 (Japan) 22:02 (India)', 'https://csc.ajis-group.co.jp/login.php', 6.376)
    ,@(1257, '10.14.2022 19:32 (Kyiv+Israel) 16:32 (UTC) 01:32 (Japan) 22:02 (India)', 'https://ru.checker-soft.com/profpoint-ru/login.php', 0.882)
    ,@(1258, '10.14.2022 19:32 (Kyiv+Israel) 16:32 (UTC) 01:32 (Japan) 22:02 (India)', 'https://www.imystery.ru/login.php', 1.085)
    ,@(1259, '10.14.2022 19:32 (Kyiv+Israel) 16:32 (UTC) 01:32 (Japan) 22:02 (India)', 'https://eu.checker-soft.com/testing/login.php', 1.095)
    ,@(1260, $null, '*****', $null)
)

$shopperData = @(
    ,@(1205, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://online.r-m.co.il/c_login.php', 0.632)
    ,@(1206, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://system.serveandcheck.com/c_login.php', 0.478)
    ,@(1207, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://eu.checker-soft.com/gfk-ukraine/c_login.php', 0.469)
    ,@(1208, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://ipsos-russia.com/c_login.php', 4.055)
    ,@(1209, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://tyaonline.com/c_login.php', 1.734)
    ,@(1210, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://marketest.checker.co.il/c_login.php', 0.278)
    ,@(1211, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://ms-online.co.il/c_login.php', 0.459)
    ,@(1212, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://www.misonline.co.il/c_login.php', 0.371)
    ,@(1213, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://tovanot.checker.co.il/c_login.php', 0.947)
    ,@(1214, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://csc.ajis-group.co.jp/c_login.php', 1.776)
    ,@(1215, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://ru.checker-soft.com/profpoint-ru/c_login.php', 0.283)
    ,@(1216, '10.13.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)', 'https://www.imystery.ru/c_login.php', 0.729)
    ,@(1217, '10.13.2022 22:03 (Kyiv+Israel) 19:03 (UTC) 04:03 (Japan) 00:33 (India)', 'https://eu.checker-soft.com/testing/c_login.php', 0.602)
    ,@(1218, $null, '*****', $null)
    ,@(1219, '10.14.2022 12:55 (Kyiv+Israel) 09:55 (UTC) 18:55 (Japan) 15:25 (India)', 'https://online.r-m.co.il/c_login.php', 0.502)
    ,@(1220, '10.14.2022 12:55 (Kyiv+Israel) 09:55 (UTC) 18:55 (Japan) 15:25 (India)', 'https://system.serveandcheck.com/c_login.php', 0.491)
    ,@(1221, '10.14.2022 12:55 (Kyiv+Israel) 09:55 (UTC) 18:55 (Japan) 15:25 (India)', 'https://eu.checker-soft.com/gfk-ukraine/c_login.php', 0.407)
    ,@(1222, '10.14.2022 12:55 (Kyiv+Israel) 09:55 (UTC) 18:55 (Japan) 15:25 (India)', 'https://ipsos-russia.com/c_login.php', 2.086)
    ,@(1223, '10.14.2022 12:55 (Kyiv+Israel) 09:55 (UTC) 18:55 (Japan) 15:25 (India)', 'https://tyaonline.com/c_login.php', 1.722)
    ,@(1224, '10.14.2022 12:56 (Kyiv+Israel) 09:56 (UTC) 18:56 (Japan) 15:26 (India)', 'https://marketest.checker.co.il/c_login.php', 0.584)
    ,@(1225, '10.14.2022 12:56 (Kyiv+Israel) 09:56 (UTC) 18:56 (Japan) 15:26 (India)', 'https://ms-online.co.il/c_login.php', 0.412)
    ,@(1226, '10.14.2022 12:56 (Kyiv+Israel) 09:56 (UTC) 18:56 (Japan) 15:26 (India)', 'https://www.misonline.co.il/c_login.php', 0.433)
    ,@(1227, '10.14.2022 12:56 (Kyiv+Israel) 09:56 (UTC) 18:56 (Japan) 15:26 (India)', 'https://tovanot.checker.co.il/c_login.php', 1.347)
    ,@(1228, '10.14.2022 12:56 (Kyiv+Israel) 09:56 (UTC) 18:56 (Japan) 15:26 (India)', 'https://csc.ajis-group.co.jp/c_login.php', 1.928)
    ,@(1229, '10.14.2022 12:56 (Kyiv+Israel) 09:56 (UTC) 18:56 (Japan) 15:26 (India)', 'https://ru.checker-soft.com/profpoint-ru/c_login.php', 0.627)
    ,@(1230, '10.14.2022 12:56 (Kyiv+Israel) 09:56 (UTC) 18:56 (Japan) 15:26 (India)', 'https://www.imystery.ru/c_login.php', 0.659)
    ,@(1231, '10.14.2022 12:56 (Kyiv+Israel) 09:56 (UTC) 18:56 (Japan) 15:26 (India)', 'https://eu.checker-soft.com/testing/c_login.php', 0.909)
    ,@(1232, $null, '*****', $null)
    ,@(1233, '10.14.2022 15:55 (Kyiv+Israel) 12:55 (UTC) 21:55 (Japan) 18:25 (India)', 'https://online.r-m.co.il/c_login.php', 0.463)
    ,@(1234, '10.14.2022 15:55 (Kyiv+Israel) 12:55 (UTC) 21:55 (Japan) 18:25 (India)', 'https://system.serveandcheck.com/c_login.php', 0.41)
    ,@(1235, '10.14.2022 15:55 (Kyiv+Israel) 12:55 (UTC) 21:55 (Japan) 18:25 (India)', 'https://eu.checker-soft.com/gfk-ukraine/c_login.php', 0.768)
    ,@(1236, '10.14.2022 15:55 (Kyiv+Israel) 12:55 (UTC) 21:55 (Japan) 18:25 (India)', 'https://ipsos-russia.com/c_login.php', 1.73)
    ,@(1237, '10.14.2022 15:56 (Kyiv+Israel) 12:56 (UTC) 21:56 (Japan) 18:26 (India)', 'https://tyaonline.com/c_login.php', 1.75)
    ,@(1238, '10.14.2022 15:56 (Kyiv+Israel) 12:56 (UTC) 21:56 (Japan) 18:26 (India)', 'https://marketest.checker.co.il/c_login.php', 0.782)
    ,@(1239, '10.14.2022 15:56 (Kyiv+Israel) 12:56 (UTC) 21:56 (Japan) 18:26 (India)', 'https://ms-online.co.il/c_login.php', 0.838)
    ,@(1240, '10.14.2022 15:56 (Kyiv+Israel) 12:56 (UTC) 21:56 (Japan) 18:26 (India)', 'https://www.misonline.co.il/c_login.php', 0.414)
    ,@(1241, '10.14.2022 15:56 (Kyiv+Israel) 12:56 (UTC) 21:56 (Japan) 18:26 (India)', 'https://tovanot.checker.co.il/c_login.php', 1.587)
    ,@(1242, '10.14.2022 15:56 (Kyiv+Israel) 12:56 (UTC) 21:56 (Japan) 18:26 (India)', 'https://csc.ajis-group.co.jp/c_login.php', 1.833)
    ,@(1243, '10.14.2022 15:57 (Kyiv+Israel) 12:57 (UTC) 21:57 (Japan) 18:27 (India)', 'https://ru.checker-soft.com/profpoint-ru/c_login.php', 0.436)
    ,@(1244, '10.14.2022 15:57 (Kyiv+Israel) 12:57 (UTC) 21:57 (Japan) 18:27 (India)', 'https://www.imystery.ru/c_login.php', 0.747)
    ,@(1245, '10.14.2022 15:57 (Kyiv+Israel) 12:57 (UTC) 21:57 (Japan) 18:27 (India)', 'https://eu.checker-soft.com/testing/c_login.php', 0.868)
    ,@(1246, $null, '*****', $null)
    ,@(1247, '10.14.2022 19:31 (Kyiv+Israel) 16:31 (UTC) 01:31 (Japan) 22:01 (India)', 'https://online.r-m.co.il/c_login.php', 0.382)
    ,@(1248, '10.14.2022 19:31 (Kyiv+Israel) 16:31 (UTC) 01:31 (Japan) 22:01 (India)', 'https://system.serveandcheck.com/c_login.php', 0.357)
    ,@(1249, '10.14.2022 19:31 (Kyiv+Israel) 16:31 (UTC) 01:31 (Japan) 22:01 (India)', 'https://eu.checker-soft.com/gfk-ukraine/c_login.php', 0.304)
    ,@(1250, '10.14.2022 19:31 (Kyiv+Israel) 16:31 (UTC) 01:31 (Japan) 22:01 (India)', 'https://ipsos-russia.com/c_login.php', 1.974)
    ,@(1251, '10.14.2022 19:31 (Kyiv+Israel) 16:31 (UTC) 01:31 (Japan) 22:01 (India)', 'https://tyaonline.com/c_login.php', 2.02)
    ,@(1252, '10.14.2022 19:32 (Kyiv+Israel) 16:32 (UTC) 01:32 (Japan) 22:02 (India)', 'https://marketest.checker.co.il/c_login.php', 0.549)
    ,@(1253, '10.14.2022 19:32 (Kyiv+Israel) 16:32 (UTC) 01:32 (Japan) 22:02 (India)', 'https://ms-online.co.il/c_login.php', 0.434)
    ,@(1254, '10.14.2022 19:32 (Kyiv+Israel) 16:32 (UTC) 01:32 (Japan) 22:02 (India)', 'https://www.misonline.co.il/c_login.php', 0.497)
    ,@(1255, '10.14.2022 19:32 (Kyiv+Israel) 16:32 (UTC) 01:32 (Japan) 22:02 (India)', 'https://tovanot.checker.co.il/c_login.php', 2.941)
    ,@(1256, '10.14.2022 19:32 (Kyiv+Israel) 16:32 (UTC) 01:32 (Japan) 22:02 (India)', 'https://csc.ajis-group.co.jp/c_login.php', 5.067)
    ,@(1257, '10.14.2022 19:32 (Kyiv+Israel) 16:32 (UTC) 01:32 (Japan) 22:02 (India)', 'https://ru.checker-soft.com/profpoint-ru/c_login.php', 0.316)
    ,@(1258, '10.14.2022 19:32 (Kyiv+Israel) 16:32 (UTC) 01:32 (Japan) 22:02 (India)', 'https://www.imystery.ru/c_login.php', 0.6889999999999999)
    ,@(1259, '10.14.2022 19:32 (Kyiv+Israel) 16:32 (UTC) 01:32 (Japan) 22:02 (India)', 'https://eu.checker-soft.com/testing/c_login.php', 0.598)
    ,@(1260, $null, '*****', $null)
)

$wsManager = $wb.Worksheets.Item("MANAGER")
foreach ($row in $managerData) {
    $r = $row[0]
    if ($row[1] -ne $null) { $wsManager.Cells.Item($r, 1).Value = $row[1] }
    if ($row[2] -ne $null) { $wsManager.Cells.Item($r, 2).Value = $row[2] }
    if ($row[3] -ne $null) { $wsManager.Cells.Item($r, 3).Value = $row[3] }
}

$wsShopper = $wb.Worksheets.Item("SHOPPER")
foreach ($row in $shopperData) {
    $r = $row[0]
    if ($row[1] -ne $null) { $wsShopper.Cells.Item($r, 1).Value = $row[1] }
    if ($row[2] -ne $null) { $wsShopper.Cells.Item($r, 2).Value = $row[2] }
    if ($row[3] -ne $null) { $wsShopper.Cells.Item($r, 3).Value = $row[3] }
}
